# Apply updated crypto price/volume data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.989.17"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.33%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.856.68"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.84%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.28%  "
$ws.Range("E6").Value = "  +0.19%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4352"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -5.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3673"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.86%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07492"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.92%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9343"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.76%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.26"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.69%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.852.85"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.15%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.686"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.69%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.401"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.64%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06886"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.13%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.005"
$ws.Range("D16").Style = "Normal"
$ws.Range("E17").Value = "  -3.02%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008999"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.83%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.003"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.79"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.41%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.988.90"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.105"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.83"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.119.28"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.60%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.011"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.88%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.22"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.73%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.29"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.357"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.38%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "112.79"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.89%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.722"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -7.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08959"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.62%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7964"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -8.30%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.811"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.28%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.020"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.86%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.167"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.53%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.004"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.22%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.117"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.42%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05412"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.44%  "
$ws.Range("E39").Value = "  -4.50%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.933"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.45%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5217"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.983"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.58%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1677"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.26%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.698"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.57%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.06698"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.42%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4854"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.23%  "
$ws.Range("E47").Value = "  -5.85%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "106.51"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.66%  "
$ws.Range("E49").Value = "  +0.14%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.905"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -7.69%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.667"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.40%  "
